$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.5422241521478246, 1.0558352295969433, 0.22523177706319525, 2.04661586751143)"
$ws.Range("C2").Value = "NIG(0.6357603174511297, 0.4407127731219181, 5.722300295976047, 4.739824404879383)"
$ws.Range("D2").Value = "F(547.8741582047332, 16.376410060635095, -9.430234757451768, 12.620450290022482)"
$ws.Range("E2").Value = "NIG(1.4032608258323256, 1.0845780052516263, 4.399606989454644, 6.06305669452199)"
